$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 2008 and 2009 data rows (rows 2 and 3). The remaining data rows
# (originally 2010..2020, rows 4..14) shift up to rows 2..12.
$ws.Rows("2:3").Delete()

# After the shift, row 12 holds the former 2020 data and row 13 is now free.
# Copy the formatting of the last existing data row onto the new row so the
# year label keeps the same cell style, then fill in the 2021 figures.
$newRow = 13
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

$ws.Cells.Item($newRow, 1).Value = "2021年"
$ws.Cells.Item($newRow, 2).Value = 1.1
$ws.Cells.Item($newRow, 5).Value = 2057.3
$ws.Cells.Item($newRow, 6).Value = 2033.6
$ws.Cells.Item($newRow, 7).Value = 2202.8
$ws.Cells.Item($newRow, 8).Value = 232.6
$ws.Cells.Item($newRow, 9).Value = 73.5
$ws.Cells.Item($newRow, 11).Value = 299.6
$ws.Cells.Item($newRow, 12).Value = 2020.9
$ws.Cells.Item($newRow, 13).Value = 287
$ws.Cells.Item($newRow, 18).Value = 43.2
$ws.Cells.Item($newRow, 20).Value = 2320.5
$ws.Cells.Item($newRow, 22).Value = 30.6
